# Update LR-pairs TPM values on the active sheet (Lpl-Sdc1.xlsx) with the
# new TPM-derived statistics, as per commit "update scripts wuth new tpm".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "G2" = 99.14059966666667
    "H2" = 297.421799
    "I2" = 0.3911422343348016
    "J2" = 0.3911422343348016
    "M2" = 0.8213140000000001
    "N2" = 2.463942
    "O2" = 0.06824749762056036
    "P2" = 0.06824749762056037
    "Q2" = 81.42556247462868
    "R2" = 732.8300622716581
    "S2" = 0.02669447870706503
    "T2" = 0.02669447870706504

    "G3" = 99.14059966666667
    "H3" = 297.421799
    "I3" = 0.3911422343348016
    "J3" = 0.3911422343348016
    "O3" = 0.2017018900182306
    "P3" = 0.2017018900182306
    "Q3" = 240.6489676477499
    "R3" = 2165.840708829749
    "S3" = 0.07889412793128313
    "T3" = 0.07889412793128314

    "G4" = 99.14059966666667
    "H4" = 297.421799
    "I4" = 0.3911422343348016
    "J4" = 0.3911422343348016
    "O4" = 0.7300506123612091
    "P4" = 0.7300506123612091
    "Q4" = 871.0177489137728
    "R4" = 7839.159740223955
    "S4" = 0.2855536276964534
    "T4" = 0.2855536276964534

    "I5" = 0.4928190063160421
    "J5" = 0.4928190063160421
    "M5" = 0.8213140000000001
    "N5" = 2.463942
    "O5" = 0.06824749762056036
    "P5" = 0.06824749762056037
    "Q5" = 102.5920017451333
    "R5" = 923.3280157062001
    "S5" = 0.033633663960921
    "T5" = 0.03363366396092101

    "I6" = 0.4928190063160421
    "J6" = 0.4928190063160421
    "O6" = 0.2017018900182306
    "P6" = 0.2017018900182306
    "S6" = 0.09940252501085203
    "T6" = 0.09940252501085203

    "I7" = 0.4928190063160421
    "J7" = 0.4928190063160421
    "O7" = 0.7300506123612091
    "P7" = 0.7300506123612091
    "S7" = 0.3597828173442691
    "T7" = 0.3597828173442691

    "G8" = 29.411685
    "H8" = 88.235055
    "I8" = 0.1160387593491562
    "J8" = 0.1160387593491562
    "M8" = 0.8213140000000001
    "N8" = 2.463942
    "O8" = 0.06824749762056036
    "P8" = 0.06824749762056037
    "Q8" = 24.15622865409
    "R8" = 217.40605788681
    "S8" = 0.007919354952574313
    "T8" = 0.007919354952574315

    "G9" = 29.411685
    "H9" = 88.235055
    "I9" = 0.1160387593491562
    "J9" = 0.1160387593491562
    "O9" = 0.2017018900182306
    "P9" = 0.2017018900182306
    "R9" = 642.5321704978051
    "S9" = 0.02340523707609543
    "T9" = 0.02340523707609543

    "G10" = 29.411685
    "H10" = 88.235055
    "I10" = 0.1160387593491562
    "J10" = 0.1160387593491562
    "O10" = 0.7300506123612091
    "P10" = 0.7300506123612091
    "S10" = 0.08471416732048645
    "T10" = 0.08471416732048645
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
